# Updates the Ccl5-Ccr3 ligand-receptor pairs table with newly recomputed
# TPM-based statistics (columns E:T, data rows 2-19). Column layout:
# E=Ligand-expressing cells, F=Ligand detection rate,
# G=Ligand average expression value, H=Ligand total expression value,
# I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate,
# M=Receptor average expression value, N=Receptor total expression value,
# O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per data row (sheet rows 2-19), values ordered E..T
$data = @(
    @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.2064979800766666, 1.85848182069, 0.01228558959668689, 0.01228558959668688),
    @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 1.795301303818, 16.157711734362, 0.1068113838833482, 0.1068113838833482),
    @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.1612008020333334, 1.4508072183, 0.009590635684199047, 0.009590635684199043),
    @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.2440064656005555, 2.196058190405, 0.01451715553921418, 0.01451715553921418),
    @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 2.121401505574333, 19.092613550169, 0.1262127031828768, 0.1262127031828768),
    @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.1904814659277778, 1.71433319335, 0.0113326876868003, 0.0113326876868003),
    @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.4407230918955555, 3.96650782706, 0.02622080385871791, 0.02622080385871791),
    @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 3.831663347065334, 34.48497012358801, 0.2279646674375945, 0.2279646674375945),
    @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.3440465415777778, 3.0964188742, 0.02046903611558213, 0.02046903611558213),
    @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.0743523505711111, 0.6691711551399999, 0.004423590314668051, 0.004423590314668051),
    @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 0.6464221677746667, 5.817799509972001, 0.0384588626800694, 0.0384588626800694),
    @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.05804249775555557, 0.5223824798000001, 0.003453236231786031, 0.003453236231786031),
    @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.3991695990394443, 3.592526391354999, 0.02374858035634029, 0.02374858035634029),
    @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 3.470395697497667, 31.233561277479, 0.2064710621466361, 0.2064710621466361),
    @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.3116081788722223, 2.80447360985, 0.01853911694038745, 0.01853911694038745),
    @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 2, 0.6666666666666666, 0.1808983333333333, 0.5426949999999999, 0.09546831801815302, 0.09546831801815302, 0.2398959127722222, 2.15906321495, 0.0142725983525257, 0.0142725983525257),
    @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 3, 1, 1.572737, 4.718211, 0.8300051930177132, 0.8300051930177132, 2.085664202723334, 18.77097782451, 0.1240865136871883, 0.1240865136871883),
    @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 1, 0.3333333333333333, 0.1412166666666667, 0.42365, 0.07452648896413371, 0.07452648896413369, 0.1872725996111111, 1.6854533965, 0.01114177630537874, 0.01114177630537874)
)

$startRow = 2
$startCol = 5  # column E

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Count) rows (E2:T19) with new TPM-based statistics"
